$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-10
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05)
$ws.Range("C2:C10").Value = 45174
